# Schedule Management .xlsx - add duration-in-hours values to the table,
# roll them up into per-sprint subtotal formulas, fix a couple of mis-styled
# cells so they match the rest of their block, and leave the selection on H11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sprint 0 (rows 7-14), subtotal in D6 ---
$ws.Range("D7").Value = 25
$ws.Range("D8").Value = 25
$ws.Range("D9").Value = 25
$ws.Range("D10").Value = 25
$ws.Range("D11").Value = 20
$ws.Range("D12").Value = 30
$ws.Range("D13").Value = 24
$ws.Range("D14").Value = 26
$ws.Range("D6").Formula = "=D7+D8+D9+D10+D11+D12+D13+D14"

# --- Sprint 1 (rows 16-20), subtotal in D15 ---
$ws.Range("D16").Value = 35
$ws.Range("D17").Value = 35
$ws.Range("D18").Value = 35
$ws.Range("D19").Value = 45
$ws.Range("D20").Value = 50
$ws.Range("D15").Formula = "=D16+D17+D18+D19+D20"

# --- Sprint 2 (rows 22-26), subtotal in D21 ---
$ws.Range("D22").Value = 35
$ws.Range("D23").Value = 35
$ws.Range("D24").Value = 35
$ws.Range("D25").Value = 45
$ws.Range("D26").Value = 50
$ws.Range("D21").Formula = "=D22+D23+D24+D25+D26"

# --- Sprint 3 (rows 28-32), subtotal in D27 ---
$ws.Range("D28").Value = 35
$ws.Range("D29").Value = 35
$ws.Range("D30").Value = 35
$ws.Range("D31").Value = 45
$ws.Range("D32").Value = 50
$ws.Range("D27").Formula = "=D28+D29+D30+D31+D32"

# --- Sprint 4 (rows 34-38), subtotal in D33 ---
$ws.Range("D34").Value = 35
$ws.Range("D35").Value = 35
$ws.Range("D36").Value = 35
$ws.Range("D37").Value = 45
$ws.Range("D38").Value = 50
$ws.Range("D33").Formula = "=D34+D35+D36+D37+D38"

# A few of the "4th row" / "5th row" cells in the later blocks had picked up
# the wrong pre-set formatting (copy/paste drift from earlier edits) - bring
# them back in line with the matching cells in the Sprint 1 block.
$ws.Range("D19").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4122) | Out-Null
$ws.Range("D31").PasteSpecial(-4122) | Out-Null
$ws.Range("D37").PasteSpecial(-4122) | Out-Null

$ws.Range("D20").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Move the selection to where the author left off.
$ws.Range("H11").Select() | Out-Null
